$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sensor window slides forward by one sample: the old first data row
# (row 2) is dropped, every remaining data row shifts up by one, and 11
# new samples are appended at the bottom (old dimension A1:C21 -> A1:C31).

# Drop the oldest sample (row 2), shifting rows 3:21 up to 2:20.
$ws.Rows.Item(2).Delete() | Out-Null

# Newly captured samples to append after the shift (these land in rows 21:31).
$newData = @(
    @(-10.10338973999023, -6.828082084655762, 5.437564849853516),
    @(9.688706398010254, -28.64854431152344, 0.4038746356964111),
    @(-21.26608657836914, 8.973570823669434, -19.81748580932617),
    @(-51.4796142578125, -10.81356239318848, -10.4013614654541),
    @(69.01158905029297, -76.07331848144531, 28.56989669799805),
    @(-14.15320301055908, 7.329947471618652, -4.596967697143555),
    @(66.06742858886719, -30.20949363708496, 25.3086986541748),
    @(-70.19232940673828, -22.23063087463379, -31.12885093688965),
    @(37.95425033569336, 1.236392974853516, -15.80910873413086),
    @(-24.9067497253418, -28.87722587585449, -10.39637756347656),
    @(-7.391507625579834, -34.15201568603516, -12.95433330535889)
)

$startRow = 21
for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newData[$i][0]
    $ws.Cells.Item($row, 2).Value = $newData[$i][1]
    $ws.Cells.Item($row, 3).Value = $newData[$i][2]
}
